$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    1  = 0.064295737209334902
    2  = -0.0099999994348607402
    3  = -0.0089999994423290985
    4  = 0.28399031917826534
    5  = -0.0059999994600161699
    6  = -0.0059999994442243576
    7  = -0.019999999356564047
    8  = -0.01999999935507546
    9  = 0.037648267317625717
    10 = -0.0059999994360069309
    11 = -0.0044999994449952396
    12 = -0.0059999994349948516
    13 = 0.017995202500526553
    14 = -0.011999999390829963
    15 = -0.0059999994264030576
    16 = -0.0059999994248971511
    17 = -0.0059999994227375453
    18 = -0.0089999994037626152
    19 = -0.0089999994396992022
    20 = -0.060862730910654861
    21 = -0.0089999994229161828
    22 = -0.0089999994224934099
    23 = -0.0089999994387888194
    24 = -0.041999999225865281
    25 = -0.041999999222049667
    26 = -0.0059999994423129976
    27 = -0.0059999994398820533
    28 = -0.0059999994286306091
    29 = -0.01199999938399543
    30 = -0.019999999330829077
    31 = -0.01499999935724361
    32 = -0.017483236092560084
    33 = -0.0059999994123254297
}

foreach ($row in $values.Keys) {
    $ws.Range("A$row").Value = $values[$row]
}
